$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "otp" header next to the existing "contact" header
$ws.Range("B1").Value = "otp"

# Add the new OTP value in column B, row 3 (row 2 of column B left blank)
$ws.Range("B3").Value = 111111

# Match the workbook's final selection/active cell state
[void]$ws.Range("B3").Select()
